$d = $word.ActiveDocument

# The "Requisitos" bullet list is the last paragraph in the document; it holds
# three runs, one per requisite, each ending in a manual line break (w:br).
# The edit moves the "LOQ4083 ..." requisite from the first position to the
# last position, leaving the other two requisites (LOB1006, LOB1019) in place
# and in their original relative order.
$requisitoText = "LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)"
$para = $d.Paragraphs.Item($d.Paragraphs.Count)

# Locate the run's range (its text plus the trailing line break) without
# disturbing anything else, then remove it from its current position.
$moveRange = $para.Range.Duplicate()
$found = $moveRange.Find.Execute($requisitoText + "^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the LOQ4083 requisito line to move"
}
$moveRange.Delete()

# Re-fetch the (now shorter) paragraph and append the requisito, as its own
# run, right before the paragraph mark -- i.e. after the last remaining
# requisito's line break.
$para = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($para.Range.End - 1, $para.Range.End - 1)
$insertionPoint.InsertAfter($requisitoText + "`v")
